# GameTestCases.xlsx - "Add files via upload" edit
#
# Fixes the leading-space typo in the "highlights grey on mouse hover"
# expected-result text (rows 2,4,6,8) and appends a second sprint's worth
# of test cases ("Select Level" feature / Easy / Medium / Hard buttons plus
# a repeat of the Exit button cases) starting at row 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the " highlights grey on mouse hover" leading-space typo -------
$ws.Range("D2").Value = "highlights grey on mouse hover"
$ws.Range("D4").Value = "highlights grey on mouse hover"
$ws.Range("D6").Value = "highlights grey on mouse hover"
$ws.Range("D8").Value = "highlights grey on mouse hover"

# --- Append the new "Select Level" sprint section (row 10 left blank) ---
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = "Select Level"
$ws.Range("C11").Value = "Easy Button Hover"
$ws.Range("D11").Value = "highlights light grey on mouse hover"
$ws.Range("E11").Value = "Pass"

$ws.Range("C12").Value = "Easy Button Press"
$ws.Range("D12").Value = "highlights grey on mouse hold"
$ws.Range("E12").Value = "Pass"

$ws.Range("C13").Value = "Medium Button Hover"
$ws.Range("D13").Value = "highlights light grey on mouse hover"
$ws.Range("E13").Value = "Pass"

$ws.Range("C14").Value = "Medium Button Press"
$ws.Range("D14").Value = "highlights grey on mouse hold"
$ws.Range("E14").Value = "Pass"

$ws.Range("C15").Value = "Hard Button Hover"
$ws.Range("D15").Value = "highlights light grey on mouse hover"
$ws.Range("E15").Value = "Pass"

$ws.Range("C16").Value = "Hard Button Press"
$ws.Range("D16").Value = "highlights grey on mouse hold"
$ws.Range("E16").Value = "Pass"

$ws.Range("C17").Value = "Exit Button Hover"
$ws.Range("D17").Value = "highlights light grey on mouse hover"
$ws.Range("E17").Value = "Pass"

$ws.Range("C18").Value = "Exit Button Press"
$ws.Range("D18").Value = "Exit application"
$ws.Range("E18").Value = "TBD"

# --- Match the author's final selection/active cell ----------------------
$ws.Range("E18").Select()
